{"js": "// Update the \"Crossing Algorithms\" bullet list: the algorithms and example\n// graphs moved from a \"Crossing-Algorithms\"/\"Crossing Graphs\" layout into a\n// shared \"Layered-Graphs/Algorithms\" + \"Layered-Graphs/Graphs\" layout, the\n// sample graph was renamed from \"1_test\" to \"crossing_test\", and two new\n// bullets were added at the end of the list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three existing bullets under the \"Crossing Algorithms\" heading\n// by their (still unique) original text so the edit is resilient to any\n// small positional differences.\nlet firstIdx = -1;\nlet secondIdx = -1;\nlet thirdIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"run both barycenter and mce\") !== -1) {\n    firstIdx = i;\n  } else if (t.indexOf(\"ditto for r_100_110_10_0_0p7\") !== -1) {\n    secondIdx = i;\n  } else if (t.indexOf(\"run sifting on\") !== -1) {\n    thirdIdx = i;\n  }\n}\n\nif (firstIdx === -1 || secondIdx === -1 || thirdIdx === -1) {\n  throw new Error(\"Could not locate the Crossing Algorithms bullets to update.\");\n}\n\nconst firstPara = paragraphs.items[firstIdx];\nconst secondPara = paragraphs.items[secondIdx];\nconst thirdPara = paragraphs.items[thirdIdx];\n\n// Rewrite the three existing bullets in place (same list formatting is kept\n// automatically since we are only replacing the paragraph's text).\nfirstPara.insertText(\n  \"run both barycenter and mce (in Layered-Graphs/Algorithms) on crossing_test (in Layered-Graphs/Graphs), need node weights for barycenter\",\n  Word.InsertLocation.replace\n);\nsecondPara.insertText(\n  \"ditto for r_100_110_10_0_0p7 (in Layered-Graphs/Graphs)\",\n  Word.InsertLocation.replace\n);\nthirdPara.insertText(\n  \"run sifting on crossing_test\",\n  Word.InsertLocation.replace\n);\n\n// Add the two new bullets after \"run sifting on crossing_test\", inheriting\n// the same ListParagraph style / numbering.\nconst fourthPara = thirdPara.insertParagraph(\"run mce on crossing_test\", Word.InsertLocation.after);\nfourthPara.insertParagraph(\"run all three on r_100_110_10_0_0p7\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Update the \"Crossing Algorithms\" bullet list: the algorithms and example\n# graphs moved from a \"Crossing-Algorithms\"/\"Crossing Graphs\" layout into a\n# shared \"Layered-Graphs/Algorithms\" + \"Layered-Graphs/Graphs\" layout, the\n# sample graph was renamed from \"1_test\" to \"crossing_test\", and two new\n# bullets were added at the end of the list.\n\n$d = $word.ActiveDocument\n\n# Locate the three existing bullets under the \"Crossing Algorithms\" heading\n# by their (still unique) original text so the edit is resilient to any\n# small positional differences.\n$firstIdx = -1\n$secondIdx = -1\n$thirdIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text\n    if ($t -like \"*run both barycenter and mce*\") {\n        $firstIdx = $i\n    } elseif ($t -like \"*ditto for r_100_110_10_0_0p7*\") {\n        $secondIdx = $i\n    } elseif ($t -like \"*run sifting on*\") {\n        $thirdIdx = $i\n    }\n}\n\nif ($firstIdx -eq -1 -or $secondIdx -eq -1 -or $thirdIdx -eq -1) {\n    throw \"Could not locate the Crossing Algorithms bullets to update.\"\n}\n\n# Rewrite the three existing bullets in place (same list formatting is kept\n# automatically since we are only replacing the paragraph's text, not the\n# paragraph mark).\n$d.Paragraphs($firstIdx).Range.Text = \"run both barycenter and mce (in Layered-Graphs/Algorithms) on crossing_test (in Layered-Graphs/Graphs), need node weights for barycenter\"\n$d.Paragraphs($secondIdx).Range.Text = \"ditto for r_100_110_10_0_0p7 (in Layered-Graphs/Graphs)\"\n$d.Paragraphs($thirdIdx).Range.Text = \"run sifting on crossing_test\"\n\n# Add the two new bullets after \"run sifting on crossing_test\", inheriting\n# the same ListParagraph style / numbering.\n$d.Paragraphs($thirdIdx).Range.InsertParagraphAfter()\n$d.Paragraphs($thirdIdx + 1).Range.Text = \"run mce on crossing_test\"\n\n$d.Paragraphs($thirdIdx + 1).Range.InsertParagraphAfter()\n$d.Paragraphs($thirdIdx + 2).Range.Text = \"run all three on r_100_110_10_0_0p7\"\n"}
